$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 14.58397638733022
$ws.Range("R2").Value = 131.255787485972
$ws.Range("S2").Value = 0.005233957458222153
$ws.Range("T2").Value = 0.005233957458222153

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 14.37421617434356
$ws.Range("R3").Value = 129.367945569092
$ws.Range("S3").Value = 0.00515867784983266
$ws.Range("T3").Value = 0.00515867784983266

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 2.145506408268778
$ws.Range("R4").Value = 19.309557674419
$ws.Range("S4").Value = 0.0007699881684516007
$ws.Range("T4").Value = 0.0007699881684516007

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 181.3673159061627
$ws.Range("R5").Value = 1632.305843155464
$ws.Range("S5").Value = 0.06508984864988313
$ws.Range("T5").Value = 0.06508984864988311

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 178.7587237223227
$ws.Range("R6").Value = 1608.828513500904
$ws.Range("S6").Value = 0.0641536663527196
$ws.Range("T6").Value = 0.06415366635271959

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 26.68166268187533
$ws.Range("R7").Value = 240.134964136878
$ws.Range("S7").Value = 0.009575624896985579
$ws.Range("T7").Value = 0.009575624896985577

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 1110.54949742002
$ws.Range("R8").Value = 9994.945476780176
$ws.Range("S8").Value = 0.3985585734900136
$ws.Range("T8").Value = 0.3985585734900135

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 1094.576549238793
$ws.Range("R9").Value = 9851.188943149136
$ws.Range("S9").Value = 0.392826136118848
$ws.Range("T9").Value = 0.392826136118848

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 163.3773259180724
$ws.Range("R10").Value = 1470.395933262652
$ws.Range("S10").Value = 0.05863352701504374
$ws.Range("T10").Value = 0.05863352701504373

Write-Output "Update complete"